# converter/input.xlsx edits:
#  - Converter date format changed to DD.MM.YYYY (column H, "Дата рождения")
#  - Document numbers (column F) updated to new sample values
#  - Client list form resize -> selection/active cell moved from K3 to H6
#  - Comments kept/restored (D1 "M or F" hint, H1 date-format hint, J1 "255 characters" hint)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Re-create the cell comments (authoring hints) -------------------------
# D1: gender hint
$ws.Range("D1").AddComment("Philipp Grigoryev:`nМ или Ж") | Out-Null
# H1: date-format hint (unchanged by this revision, but must stay present)
$ws.Range("H1").AddComment("Philipp Grigoryev:`nв формате ДД.ММ.ГГГГ") | Out-Null
# J1: max length hint
$ws.Range("J1").AddComment("255 символов") | Out-Null

# --- Row 2: Гальмутдинов -----------------------------------------------------
$ws.Range("F2").Value = 4507120315
$ws.Range("H2").Value = "30.02.1960"

# --- Row 3: Литров -----------------------------------------------------------
$ws.Range("F3").Value = 4004047028
# "01.01.1970" is a real calendar date, so typing it directly would be
# auto-recognised by Excel and stored as a date serial number. Build it as a
# text formula first, then collapse the formula down to its literal text
# value (copy / paste-values) so the cell keeps storing plain text, exactly
# like the other text cells in the sheet.
$ws.Range("H3").Formula = "=""01.01.1970"""
$ws.Range("H3").Copy()
$ws.Range("H3").PasteSpecial(-4163)

# --- Row 4: Рюмочкина ----------------------------------------------------
$ws.Range("F4").Value = 4507887863
$ws.Range("H4").Formula = "=""02.02.1970"""
$ws.Range("H4").Copy()
$ws.Range("H4").PasteSpecial(-4163)

# --- Row 5: Бухайло -------------------------------------------------------
# Document number becomes a genuine number (was stored as text before).
$ws.Range("F5").Value = 4507887843
$ws.Range("H5").Formula = "=""03.03.1950"""
$ws.Range("H5").Copy()
$ws.Range("H5").PasteSpecial(-4163)

# --- Client list form resize: active selection moved -----------------------
$ws.Range("H6").Select() | Out-Null
